# Convert all v1 to v2
#
# v1 layout: library_content | controls | implementation_groups
# v2 layout: library_meta | controls_meta | controls_content |
#            implementation_groups_meta | implementation_groups_content
#
# NOTE: worksheet object references returned by Worksheets.Item(...) in
# this host are resolved *positionally* and go stale across any
# structural operation (Add/Move) that shifts tab order. So every sheet
# is re-fetched by (stable) name immediately before it is used, and all
# sheet-Add operations happen before any other content work touches
# sheets whose position could shift.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Rename the three original sheets to their v2 "content"/final names.
#    (Renaming does not change tab order, so this is safe to do first.)
# ---------------------------------------------------------------------
$wb.Worksheets.Item("library_content").Name = "library_meta"
$wb.Worksheets.Item("controls").Name = "controls_content"
$wb.Worksheets.Item("implementation_groups").Name = "implementation_groups_content"

# ---------------------------------------------------------------------
# 2) Insert the two brand-new "meta" sheets right before their matching
#    "content" sheet, re-fetching by name each time.
# ---------------------------------------------------------------------
$controlsMeta = $wb.Worksheets.Add($wb.Worksheets.Item("controls_content"))
$controlsMeta.Name = "controls_meta"

$igMeta = $wb.Worksheets.Add($wb.Worksheets.Item("implementation_groups_content"))
$igMeta.Name = "implementation_groups_meta"

# Final expected order: library_meta, controls_meta, controls_content,
# implementation_groups_meta, implementation_groups_content

# ---------------------------------------------------------------------
# 3) library_meta: trim to the pure "library" record (10 rows x 2 cols)
# ---------------------------------------------------------------------
$libraryDescription = @"
AirCyber is the AeroSpace and Defense official standard for Cybersecurity maturity evaluation and increase built by Airbus, Dassault Aviation, Safran and Thales to help the AeroSpace SupplyChain to be more resilient.
Their joint venture BoostAeroSpace is offering this extract of the AirCyber maturity level matrix to provide further details on this standard, the questions and the AirCyber maturity levels they are associated to.
AirCyber program uses this maturity level matrix as the base of the cyber maturity evaluation as is the evaluation activity is the very starting point for any cyber maturity progression. Being aware of the problems is the mandatory very first knowledge a company shall know to decide to launch a cybersecurity company program.
Source: https://boostaerospace.com/aircyber/
"@

$libraryCopyright = @"
© Boost Aerospace
This work is licensed under a Creative Commons Attribution-NonCommercial-ShareAlike 4.0 International License. Any commercial use of this work must be contracted with BoostAeroSpace.
Permission given to include AirCyber in CISO Assistant.
"@

$wsLibraryMeta = $wb.Worksheets.Item("library_meta")

# Clear everything first (old sheet had rows 1-15, cols A-C; new one only
# needs rows 1-10, cols A-B) so no stray cells survive.
$wsLibraryMeta.Cells.Clear()

$wsLibraryMeta.Range("A1").Value = "type"
$wsLibraryMeta.Range("B1").Value = "library"

$wsLibraryMeta.Range("A2").Value = "urn"
$wsLibraryMeta.Range("B2").Value = "urn:intuitem:risk:library:aircyber-v1.5.2"

$wsLibraryMeta.Range("A3").Value = "version"
$wsLibraryMeta.Range("B3").NumberFormat = "@"
$wsLibraryMeta.Range("B3").Value = "1"

$wsLibraryMeta.Range("A4").Value = "locale"
$wsLibraryMeta.Range("B4").Value = "en"

$wsLibraryMeta.Range("A5").Value = "ref_id"
$wsLibraryMeta.Range("B5").Value = "AirCyber-v1.5.2"

$wsLibraryMeta.Range("A6").Value = "name"
$wsLibraryMeta.Range("B6").Value = "Public AirCyber Maturity Level Matrix"

$wsLibraryMeta.Range("A7").Value = "description"
$wsLibraryMeta.Range("B7").Value = $libraryDescription

$wsLibraryMeta.Range("A8").Value = "copyright"
$wsLibraryMeta.Range("B8").Value = $libraryCopyright

$wsLibraryMeta.Range("A9").Value = "provider"
$wsLibraryMeta.Range("B9").Value = "Boost Aerospace"

$wsLibraryMeta.Range("A10").Value = "packager"
$wsLibraryMeta.Range("B10").Value = "intuitem"

# ---------------------------------------------------------------------
# 4) controls_meta: brand-new "framework" record (7 rows x 2 cols)
# ---------------------------------------------------------------------
$frameworkDescription = @"
AirCyber is the AeroSpace and Defense official standard for Cybersecurity maturity evaluation and increase built by Airbus, Dassault Aviation, Safran and Thales to help the AeroSpace SupplyChain to be more resilient.
Their joint venture BoostAeroSpace is offering this extract of the AirCyber maturity level matrix to provide further details on this standard, the questions and the AirCyber maturity levels they are associated to.
AirCyber program uses this maturity level matrix as the base of the cyber maturity evaluation as is the evaluation activity is the very starting point for any cyber maturity progression. Being aware of the problems is the mandatory very first knowledge a company shall know to decide to launch a cybersecurity company program.
Source: https://boostaerospace.com/aircyber/
"@

$wsControlsMeta = $wb.Worksheets.Item("controls_meta")

$wsControlsMeta.Range("A1").Value = "type"
$wsControlsMeta.Range("B1").Value = "framework"

$wsControlsMeta.Range("A2").Value = "base_urn"
$wsControlsMeta.Range("B2").Value = "urn:intuitem:risk:req_node:aircyber-v1.5.2"

$wsControlsMeta.Range("A3").Value = "urn"
$wsControlsMeta.Range("B3").Value = "urn:intuitem:risk:framework:aircyber-v1.5.2"

$wsControlsMeta.Range("A4").Value = "ref_id"
$wsControlsMeta.Range("B4").Value = "AirCyber-v1.5.2"

$wsControlsMeta.Range("A5").Value = "name"
$wsControlsMeta.Range("B5").Value = "Public AirCyber Maturity Level Matrix"

$wsControlsMeta.Range("A6").Value = "description"
$wsControlsMeta.Range("B6").Value = $frameworkDescription

$wsControlsMeta.Range("A7").Value = "implementation_groups_definition"
$wsControlsMeta.Range("B7").Value = "implementation_groups"

# ---------------------------------------------------------------------
# 5) implementation_groups_meta: brand-new tiny record (2 rows x 2 cols)
# ---------------------------------------------------------------------
$wsIGMeta = $wb.Worksheets.Item("implementation_groups_meta")

$wsIGMeta.Range("A1").Value = "type"
$wsIGMeta.Range("B1").Value = "implementation_groups"

$wsIGMeta.Range("A2").Value = "name"
$wsIGMeta.Range("B2").Value = "implementation_groups"

# ---------------------------------------------------------------------
# 6) implementation_groups_content: strip the stray empty B/C cells on
#    rows 2-4 (Bronze/Silver/Gold), leaving just column A populated
#    under the unchanged ref_id/name/description header row.
# ---------------------------------------------------------------------
$wsIGContent = $wb.Worksheets.Item("implementation_groups_content")
$wsIGContent.Range("B2:C4").Clear()

# controls_content keeps its original data untouched (only the rename in
# step 1 applies to it).
